# Rename scenery to scenario
#
# 1. The "scenery_array" sheet becomes "scenario_array".
# 2. Its companion _xlnm.Sheet_Title defined name is updated to match.
# 3. The "Demand@lmax" sheet's header cell A1 ("scenery") becomes "scenario"
#    and its row is given an explicit height.
# 4. The active/selected tab moves from stage_array to Demand@lmax, with the
#    selection on that sheet reset to A1.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the sheet and fix up its Sheet_Title defined name ---
$scenerySheet = $wb.Worksheets.Item("scenery_array")
$scenerySheet.Name = "scenario_array"

$titleName = $wb.Names.Item("scenario_array!Sheet_Title")
$titleName.RefersTo = '="scenario_array"'

# --- 3: update the header label + row height on Demand@lmax ---
$lmaxSheet = $wb.Worksheets.Item("Demand@lmax")
$lmaxSheet.Range("A1").Value = "scenario"
$lmaxSheet.Rows.Item(1).RowHeight = 15

# --- 4: make Demand@lmax the active sheet with A1 selected ---
$lmaxSheet.Activate() | Out-Null
$lmaxSheet.Range("A1").Select() | Out-Null
